$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# Extend the table with 4 more time-registration rows (31-34), copying the
# formatting (styles) of the last existing row (30) so no new number formats
# get introduced.
$ws.Range("A30:I30").Copy()
$ws.Range("A31:I34").PasteSpecial(-4122)
$ws.Range("C31:C34").Clear()
$ws.Range("E31:E34").Clear()

$data = @(
    @{ Row = 31; Date = 42817; Aktivitet = "Implementeret ReferencespaendingTest"; Start = 0.33333333333333331; Slut = 0.42708333333333331; Tid = 0.09375 },
    @{ Row = 32; Date = 42817; Aktivitet = "Aflevering af projektet";              Start = 0.4375;              Slut = 0.4861111111111111;  Tid = 0.04861111111111111 },
    @{ Row = 33; Date = 42817; Aktivitet = "Undervisning";                         Start = 0.52083333333333337; Slut = 0.58333333333333337; Tid = 0.0625 },
    @{ Row = 34; Date = 42817; Aktivitet = "sigmaRef regner forkert";              Start = 0.58333333333333337; Slut = 0.64583333333333337; Tid = 0.0625 }
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "AS"
    $ws.Cells.Item($row, 4).Value = "?"
    $ws.Cells.Item($row, 6).Value = $r.Aktivitet
    $ws.Cells.Item($row, 7).Value = $r.Start
    $ws.Cells.Item($row, 8).Value = $r.Slut
    $ws.Cells.Item($row, 9).Value = $r.Tid
}

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:I34"))

$ws.Range("A13").Select()
$ws.Range("E30").Select()
